# Auto-generated edit script: updates Leve profit-tracking values across all 8 sheets
# to match the scheduled-runner refresh described in the commit message.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 12822511   # H100: was 13335391
$ws.Cells.Item(100, 9).Value = 22223536   # I100: was 23810896
$ws.Cells.Item(100, 11).Value = 22223536   # K100: was 23810896
$ws.Cells.Item(100, 13).Value = -22222995   # M100: was -23810355
$ws.Cells.Item(137, 8).Value = 1599.238   # H137: was 1406.4138
$ws.Cells.Item(137, 9).Value = 1253.4445   # I137: was 1029.4615
$ws.Cells.Item(137, 10).Value = 1858.5834   # J137: was 1712.6875
$ws.Cells.Item(137, 11).Value = 3760.3335   # K137: was 3088.3845
$ws.Cells.Item(137, 12).Value = 5575.7502   # L137: was 5138.0625
$ws.Cells.Item(137, 13).Value = -1210.3335   # M137: was -538.3844999999997
$ws.Cells.Item(137, 14).Value = -10675.7502   # N137: was -10238.0625

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 569748.0600000001   # H32: was 548270.5
$ws.Cells.Item(32, 9).Value = 5264.3555   # I32: was 5156.5435
$ws.Cells.Item(32, 10).Value = 4803376   # J32: was 4117305.2
$ws.Cells.Item(32, 11).Value = 5264.3555   # K32: was 5156.5435
$ws.Cells.Item(32, 12).Value = 4803376   # L32: was 4117305.2
$ws.Cells.Item(32, 13).Value = -4977.3555   # M32: was -4869.5435
$ws.Cells.Item(32, 14).Value = -4803950   # N32: was -4117879.2
$ws.Cells.Item(61, 8).Value = 1767.125   # H61: was 1694.6666
$ws.Cells.Item(61, 9).Value = 1824.6   # I61: was 1747.6
$ws.Cells.Item(61, 10).Value = 1671.3334   # J61: was 1628.5
$ws.Cells.Item(61, 11).Value = 1824.6   # K61: was 1747.6
$ws.Cells.Item(61, 12).Value = 1671.3334   # L61: was 1628.5
$ws.Cells.Item(61, 13).Value = -1612.6   # M61: was -1535.6
$ws.Cells.Item(61, 14).Value = -2095.3334   # N61: was -2052.5
$ws.Cells.Item(74, 8).Value = 899.04254   # H74: was 902.44684
$ws.Cells.Item(74, 9).Value = 716.2121   # I74: was 743.0645
$ws.Cells.Item(74, 10).Value = 1330   # J74: was 1211.25
$ws.Cells.Item(74, 11).Value = 716.2121   # K74: was 743.0645
$ws.Cells.Item(74, 12).Value = 1330   # L74: was 1211.25
$ws.Cells.Item(74, 13).Value = 157.7879   # M74: was 130.9355
$ws.Cells.Item(74, 14).Value = -3078   # N74: was -2959.25
$ws.Cells.Item(77, 8).Value = 899.04254   # H77: was 902.44684
$ws.Cells.Item(77, 9).Value = 716.2121   # I77: was 743.0645
$ws.Cells.Item(77, 10).Value = 1330   # J77: was 1211.25
$ws.Cells.Item(77, 11).Value = 3581.0605   # K77: was 3715.3225
$ws.Cells.Item(77, 12).Value = 6650   # L77: was 6056.25
$ws.Cells.Item(77, 13).Value = 786.9395000000004   # M77: was 652.6775000000002
$ws.Cells.Item(77, 14).Value = -15386   # N77: was -14792.25
$ws.Cells.Item(82, 8).Value = 18000   # H82: was 15258
$ws.Cells.Item(82, 9).Value = 0   # I82: was 2164
$ws.Cells.Item(82, 10).Value = 18000   # J82: was 21805
$ws.Cells.Item(82, 11).Value = 0   # K82: was 2164
$ws.Cells.Item(82, 12).Value = 18000   # L82: was 21805
$ws.Cells.Item(82, 13).ClearContents()   # M82: was -1803
$ws.Cells.Item(82, 14).Value = -18722   # N82: was -22527
$ws.Cells.Item(85, 8).Value = 18000   # H85: was 15258
$ws.Cells.Item(85, 9).Value = 0   # I85: was 2164
$ws.Cells.Item(85, 10).Value = 18000   # J85: was 21805
$ws.Cells.Item(85, 11).Value = 0   # K85: was 2164
$ws.Cells.Item(85, 12).Value = 18000   # L85: was 21805
$ws.Cells.Item(85, 13).ClearContents()   # M85: was -916
$ws.Cells.Item(85, 14).Value = -20496   # N85: was -24301
$ws.Cells.Item(132, 8).Value = 29443048   # H132: was 30335224
$ws.Cells.Item(132, 9).Value = 33334958   # I132: was 34484396
$ws.Cells.Item(132, 11).Value = 100004874   # K132: was 103453188
$ws.Cells.Item(132, 13).Value = -100002344   # M132: was -103450658
$ws.Cells.Item(136, 8).Value = 1767.125   # H136: was 1694.6666
$ws.Cells.Item(136, 9).Value = 1824.6   # I136: was 1747.6
$ws.Cells.Item(136, 10).Value = 1671.3334   # J136: was 1628.5
$ws.Cells.Item(136, 11).Value = 5473.799999999999   # K136: was 5242.799999999999
$ws.Cells.Item(136, 12).Value = 5014.0002   # L136: was 4885.5
$ws.Cells.Item(136, 13).Value = -2923.799999999999   # M136: was -2692.799999999999
$ws.Cells.Item(136, 14).Value = -10114.0002   # N136: was -9985.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(88, 8).Value = 0   # H88: was 15666.667
$ws.Cells.Item(88, 10).Value = 0   # J88: was 15666.667
$ws.Cells.Item(88, 12).Value = 0   # L88: was 15666.667
$ws.Cells.Item(88, 14).ClearContents()   # N88: was -16478.667
$ws.Cells.Item(91, 8).Value = 0   # H91: was 15666.667
$ws.Cells.Item(91, 10).Value = 0   # J91: was 15666.667
$ws.Cells.Item(91, 12).Value = 0   # L91: was 15666.667
$ws.Cells.Item(91, 14).ClearContents()   # N91: was -18474.667
$ws.Cells.Item(134, 8).Value = 10084.95   # H134: was 9260.227999999999
$ws.Cells.Item(134, 9).Value = 3392.3333   # I134: was 3212.3125
$ws.Cells.Item(134, 10).Value = 30162.8   # J134: was 25388
$ws.Cells.Item(134, 11).Value = 10176.9999   # K134: was 9636.9375
$ws.Cells.Item(134, 12).Value = 90488.39999999999   # L134: was 76164
$ws.Cells.Item(134, 13).Value = -7641.999899999999   # M134: was -7101.9375
$ws.Cells.Item(134, 14).Value = -95558.39999999999   # N134: was -81234

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3089.524   # H31: was 2787.1072
$ws.Cells.Item(31, 9).Value = 2881.1765   # I31: was 2349.5217
$ws.Cells.Item(31, 10).Value = 3975   # J31: was 4800
$ws.Cells.Item(31, 11).Value = 2881.1765   # K31: was 2349.5217
$ws.Cells.Item(31, 12).Value = 3975   # L31: was 4800
$ws.Cells.Item(31, 13).Value = -2586.1765   # M31: was -2054.5217
$ws.Cells.Item(31, 14).Value = -4565   # N31: was -5390
$ws.Cells.Item(34, 8).Value = 3089.524   # H34: was 2787.1072
$ws.Cells.Item(34, 9).Value = 2881.1765   # I34: was 2349.5217
$ws.Cells.Item(34, 10).Value = 3975   # J34: was 4800
$ws.Cells.Item(34, 11).Value = 2881.1765   # K34: was 2349.5217
$ws.Cells.Item(34, 12).Value = 3975   # L34: was 4800
$ws.Cells.Item(34, 13).Value = -2679.1765   # M34: was -2147.5217
$ws.Cells.Item(34, 14).Value = -4379   # N34: was -5204
$ws.Cells.Item(58, 8).Value = 1376.1034   # H58: was 1674.6
$ws.Cells.Item(58, 9).Value = 1211.0952   # I58: was 1409.9412
$ws.Cells.Item(58, 10).Value = 1809.25   # J58: was 2237
$ws.Cells.Item(58, 11).Value = 1211.0952   # K58: was 1409.9412
$ws.Cells.Item(58, 12).Value = 1809.25   # L58: was 2237
$ws.Cells.Item(58, 13).Value = -1008.0952   # M58: was -1206.9412
$ws.Cells.Item(58, 14).Value = -2215.25   # N58: was -2643
$ws.Cells.Item(110, 8).Value = 40000   # H110: was 0
$ws.Cells.Item(110, 10).Value = 40000   # J110: was 0
$ws.Cells.Item(110, 12).Value = 40000   # L110: was 0
$ws.Cells.Item(110, 14).Value = -48180   # N110: was (empty)
$ws.Cells.Item(132, 8).Value = 65874.25   # H132: was 128998.375
$ws.Cells.Item(132, 9).Value = 2999.1   # I132: was 3663.3333
$ws.Cells.Item(132, 10).Value = 170666.17   # J132: was 204199.4
$ws.Cells.Item(132, 11).Value = 8997.299999999999   # K132: was 10989.9999
$ws.Cells.Item(132, 12).Value = 511998.51   # L132: was 612598.2
$ws.Cells.Item(132, 13).Value = -6467.299999999999   # M132: was -8459.999899999999
$ws.Cells.Item(132, 14).Value = -517058.51   # N132: was -617658.2
$ws.Cells.Item(134, 8).Value = 8013.5   # H134: was 7622.72
$ws.Cells.Item(134, 9).Value = 1754.091   # I134: was 1024
$ws.Cells.Item(134, 10).Value = 13309.923   # J134: was 28518.666
$ws.Cells.Item(134, 11).Value = 5262.272999999999   # K134: was 3072
$ws.Cells.Item(134, 12).Value = 39929.769   # L134: was 85555.99800000001
$ws.Cells.Item(134, 13).Value = -2727.272999999999   # M134: was -537
$ws.Cells.Item(134, 14).Value = -44999.769   # N134: was -90625.99800000001
$ws.Cells.Item(136, 8).Value = 1376.1034   # H136: was 1674.6
$ws.Cells.Item(136, 9).Value = 1211.0952   # I136: was 1409.9412
$ws.Cells.Item(136, 10).Value = 1809.25   # J136: was 2237
$ws.Cells.Item(136, 11).Value = 3633.2856   # K136: was 4229.8236
$ws.Cells.Item(136, 12).Value = 5427.75   # L136: was 6711
$ws.Cells.Item(136, 13).Value = -1083.2856   # M136: was -1679.8236
$ws.Cells.Item(136, 14).Value = -10527.75   # N136: was -11811

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(99, 8).Value = 1000   # H99: was 1349.75
$ws.Cells.Item(99, 9).Value = 1000   # I99: was 800
$ws.Cells.Item(99, 10).Value = 0   # J99: was 1533
$ws.Cells.Item(99, 11).Value = 3000   # K99: was 2400
$ws.Cells.Item(99, 12).Value = 0   # L99: was 4599
$ws.Cells.Item(99, 13).Value = -754   # M99: was -154
$ws.Cells.Item(99, 14).ClearContents()   # N99: was -9091
$ws.Cells.Item(113, 8).Value = 917.57574   # H113: was 841.55554
$ws.Cells.Item(113, 9).Value = 642.0833   # I113: was 636.7143
$ws.Cells.Item(113, 10).Value = 955.5747   # J113: was 875.2941
$ws.Cells.Item(113, 11).Value = 1926.2499   # K113: was 1910.1429
$ws.Cells.Item(113, 12).Value = 2866.7241   # L113: was 2625.8823
$ws.Cells.Item(113, 13).Value = 243.7501   # M113: was 259.8571000000002
$ws.Cells.Item(113, 14).Value = -7206.724099999999   # N113: was -6965.882299999999

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1180.1936   # H102: was 2086
$ws.Cells.Item(102, 9).Value = 951.2759   # I102: was 1655
$ws.Cells.Item(102, 10).Value = 4499.5   # J102: was 3666.3333
$ws.Cells.Item(102, 11).Value = 951.2759   # K102: was 1655
$ws.Cells.Item(102, 12).Value = 4499.5   # L102: was 3666.3333
$ws.Cells.Item(102, 13).Value = 670.7241   # M102: was -33
$ws.Cells.Item(102, 14).Value = -7743.5   # N102: was -6910.3333
$ws.Cells.Item(104, 8).Value = 27835.5   # H104: was 0
$ws.Cells.Item(104, 10).Value = 27835.5   # J104: was 0
$ws.Cells.Item(104, 12).Value = 27835.5   # L104: was 0
$ws.Cells.Item(104, 14).Value = -34823.5   # N104: was (empty)
$ws.Cells.Item(132, 8).Value = 61029.293   # H132: was 49667.145
$ws.Cells.Item(132, 9).Value = 2013.7273   # I132: was 1897.3572
$ws.Cells.Item(132, 10).Value = 169224.5   # J132: was 145206.72
$ws.Cells.Item(132, 11).Value = 6041.1819   # K132: was 5692.071599999999
$ws.Cells.Item(132, 12).Value = 507673.5   # L132: was 435620.16
$ws.Cells.Item(132, 13).Value = -3511.1819   # M132: was -3162.071599999999
$ws.Cells.Item(132, 14).Value = -512733.5   # N132: was -440680.16

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2071   # H40: was 2083.963
$ws.Cells.Item(40, 9).Value = 1968.7273   # I40: was 1976.762
$ws.Cells.Item(40, 10).Value = 2521   # J40: was 2459.1667
$ws.Cells.Item(40, 11).Value = 1968.7273   # K40: was 1976.762
$ws.Cells.Item(40, 12).Value = 2521   # L40: was 2459.1667
$ws.Cells.Item(40, 13).Value = -1832.7273   # M40: was -1840.762
$ws.Cells.Item(40, 14).Value = -2793   # N40: was -2731.1667
$ws.Cells.Item(61, 8).Value = 1540.2354   # H61: was 1445.3158
$ws.Cells.Item(61, 9).Value = 908.6667   # I61: was 877.2308
$ws.Cells.Item(61, 10).Value = 3056   # J61: was 2676.1667
$ws.Cells.Item(61, 11).Value = 908.6667   # K61: was 877.2308
$ws.Cells.Item(61, 12).Value = 3056   # L61: was 2676.1667
$ws.Cells.Item(61, 13).Value = -706.6667   # M61: was -675.2308
$ws.Cells.Item(61, 14).Value = -3460   # N61: was -3080.1667
$ws.Cells.Item(113, 8).Value = 1540.2354   # H113: was 1445.3158
$ws.Cells.Item(113, 9).Value = 908.6667   # I113: was 877.2308
$ws.Cells.Item(113, 10).Value = 3056   # J113: was 2676.1667
$ws.Cells.Item(113, 11).Value = 908.6667   # K113: was 877.2308
$ws.Cells.Item(113, 12).Value = 3056   # L113: was 2676.1667
$ws.Cells.Item(113, 13).Value = 1261.3333   # M113: was 1292.7692
$ws.Cells.Item(113, 14).Value = -7396   # N113: was -7016.1667
$ws.Cells.Item(132, 8).Value = 628002.9399999999   # H132: was 387090
$ws.Cells.Item(132, 9).Value = 835545.3   # I132: was 418847.6
$ws.Cells.Item(132, 10).Value = 5375.75   # J132: was 5999
$ws.Cells.Item(132, 11).Value = 2506635.9   # K132: was 1256542.8
$ws.Cells.Item(132, 12).Value = 16127.25   # L132: was 17997
$ws.Cells.Item(132, 13).Value = -2504105.9   # M132: was -1254012.8
$ws.Cells.Item(132, 14).Value = -21187.25   # N132: was -23057
$ws.Cells.Item(136, 8).Value = 71440760   # H136: was 32263958
$ws.Cells.Item(136, 9).Value = 21678.8   # I136: was 7685.8
$ws.Cells.Item(136, 10).Value = 111118024   # J136: was 62504216
$ws.Cells.Item(136, 11).Value = 65036.39999999999   # K136: was 23057.4
$ws.Cells.Item(136, 12).Value = 333354072   # L136: was 187512648
$ws.Cells.Item(136, 13).Value = -62486.39999999999   # M136: was -20507.4
$ws.Cells.Item(136, 14).Value = -333359172   # N136: was -187517748
$ws.Cells.Item(140, 8).Value = 100246620   # H140: was 111379020
$ws.Cells.Item(140, 10).Value = 298652.25   # J140: was 333459.84
$ws.Cells.Item(140, 12).Value = 298652.25   # L140: was 333459.84
$ws.Cells.Item(140, 14).Value = -309012.25   # N140: was -343819.84

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 44404676   # H132: was 45275372
$ws.Cells.Item(132, 9).Value = 107144264   # I132: was 97827560
$ws.Cells.Item(132, 10).Value = 1903663.4   # J132: was 2107502.2
$ws.Cells.Item(132, 11).Value = 321432792   # K132: was 293482680
$ws.Cells.Item(132, 12).Value = 5710990.199999999   # L132: was 6322506.600000001
$ws.Cells.Item(132, 13).Value = -321430262   # M132: was -293480150
$ws.Cells.Item(132, 14).Value = -5716050.199999999   # N132: was -6327566.600000001
$ws.Cells.Item(133, 8).Value = 44595   # H133: was 31808
$ws.Cells.Item(133, 10).Value = 49190   # J133: was 29760
$ws.Cells.Item(133, 12).Value = 49190   # L133: was 29760
$ws.Cells.Item(133, 14).Value = -59310   # N133: was -39880
